# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-cell updates (Price / Volume columns) ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.476.29'
$ws.Range("E2").Value = '  -0.58%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.902.62'
$ws.Range("E3").Value = '  -0.75%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.84'
$ws.Range("E5").Value = '  -2.54%  '

$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4789'
$ws.Range("E7").Value = '  +2.46%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4066'
$ws.Range("E8").Value = '  -1.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08072'
$ws.Range("E9").Value = '  +0.54%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.003'
$ws.Range("E10").Value = '  -1.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.34'
$ws.Range("E11").Value = '  +4.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.897.15'
$ws.Range("E12").Value = '  +0.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.949'
$ws.Range("E13").Value = '  -0.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.072'
$ws.Range("E14").Value = '  -1.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.84'
$ws.Range("E15").Value = '  -0.22%  '

$ws.Range("E16").Value = '  +0.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06697'
$ws.Range("E17").Value = '  +1.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001032'
$ws.Range("E18").Value = '  -0.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.62'
$ws.Range("E19").Value = '  -1.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.480.20'
$ws.Range("E21").Value = '  -0.40%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.540'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.77'
$ws.Range("E23").Value = '  +1.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.167'
$ws.Range("E24").Value = '  -2.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.121.59'
$ws.Range("E25").Value = '  +0.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.70'
$ws.Range("E26").Value = '  -0.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.84'
$ws.Range("E27").Value = '  -0.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.085'
$ws.Range("E28").Value = '  +5.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.095'
$ws.Range("E29").Value = '  -2.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.32'
$ws.Range("E30").Value = '  +0.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.036'
$ws.Range("E31").Value = '  -3.49%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09511'
$ws.Range("E32").Value = '  +0.44%  '

$ws.Range("E38").Value = '  -0.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5878'
$ws.Range("E39").Value = '  -0.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.917'
$ws.Range("E40").Value = '  -6.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1845'
$ws.Range("E41").Value = '  -0.04%  '

$ws.Range("E42").Value = '  -0.27%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.418'
$ws.Range("E43").Value = '  +2.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.280'
$ws.Range("E44").Value = '  +0.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07796'
$ws.Range("E45").Value = '  +3.62%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.24'
$ws.Range("E46").Value = '  +0.41%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5530'
$ws.Range("E47").Value = '  -0.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.922'
$ws.Range("E48").Value = '  -0.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '113.68'
$ws.Range("E49").Value = '  +0.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.2936'
$ws.Range("E50").Value = '  -2.10%  '

# --- Full row updates (coin reordering / replacement) ---
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.432'
$ws.Range("E33").Value = '  +0.45%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.393'
$ws.Range("E34").Value = '  -2.73%  '

$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.544'
$ws.Range("E35").Value = '  -0.88%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06075'
$ws.Range("E36").Value = '  -0.91%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02251'
$ws.Range("E37").Value = '  -0.80%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.17'
$ws.Range("E51").Value = '  +0.74%  '
